$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp footer text
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 14:46"

# Row 16
$ws.Range("A16").Value = "A Coruña"
$ws.Range("B16").Value = 384
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 326
$ws.Range("E16").Value = 3

# Row 17
$ws.Range("A17").Value = "Alacant/Alicante"
$ws.Range("B17").Value = 372
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 348
$ws.Range("E17").Value = 17

# Row 18
$ws.Range("A18").Value = "Pontevedra"
$ws.Range("B18").Value = 348
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 262
$ws.Range("E18").Value = 2

# Row 19
$ws.Range("A19").Value = "Zaragoza"
$ws.Range("B19").Value = 329
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 315
$ws.Range("E19").Value = 14

# Row 20
$ws.Range("A20").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B20").Value = 325
$ws.Range("C20").Value = 21
$ws.Range("D20").Value = 312
$ws.Range("E20").Value = 13

# Row 21
$ws.Range("A21").Value = "Murcia"
$ws.Range("B21").Value = 296
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 213
$ws.Range("E21").Value = 1

# Row 22
$ws.Range("A22").Value = "Granada"
$ws.Range("B22").Value = 289
$ws.Range("C22").Value = 72
$ws.Range("D22").Value = 276
$ws.Range("E22").Value = 13

# Row 23
$ws.Range("A23").Value = "Cantabria"
$ws.Range("B23").Value = 282
$ws.Range("C23").Value = 11
$ws.Range("D23").Value = 200
$ws.Range("E23").Value = 5

# Row 24
$ws.Range("A24").Value = "Burgos"
$ws.Range("B24").Value = 269
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 175
$ws.Range("E24").Value = 16

# Row 25
$ws.Range("A25").Value = "Salamanca"
$ws.Range("B25").Value = 265
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 180
$ws.Range("E25").Value = 21

# Row 27
$ws.Range("A27").Value = "Tenerife"
$ws.Range("B27").Value = 262
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 211
$ws.Range("E27").Value = 5

# Row 28
$ws.Range("A28").Value = "Sevilla"
$ws.Range("B28").Value = 245
$ws.Range("C28").Value = 72
$ws.Range("D28").Value = 243
$ws.Range("E28").Value = 2

# Row 29
$ws.Range("A29").Value = "Caceres"
$ws.Range("B29").Value = 243
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 231
$ws.Range("E29").Value = 12

# Row 30
$ws.Range("A30").Value = "Valladolid"
$ws.Range("B30").Value = 241
$ws.Range("C30").Value = 13
$ws.Range("D30").Value = 193
$ws.Range("E30").Value = 11

# Row 39
$ws.Range("A39").Value = "Gran Canaria"
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 96
$ws.Range("E39").Value = 2

# Row 40
$ws.Range("A40").Value = "Soria"
$ws.Range("B40").Value = 119
$ws.Range("C40").Value = 5
$ws.Range("D40").Value = 71
$ws.Range("E40").Value = 8

# Row 41
$ws.Range("A41").Value = "Avila"
$ws.Range("B41").Value = 114
$ws.Range("C41").Value = 14
$ws.Range("D41").Value = 55
$ws.Range("E41").Value = 9

# Row 42
$ws.Range("A42").Value = "Badajoz"
$ws.Range("B42").Value = 111
$ws.Range("C42").Value = 5
$ws.Range("D42").Value = 104
$ws.Range("E42").Value = 2

# Row 43
$ws.Range("A43").Value = "Castello/Castellon"
$ws.Range("B43").Value = 104
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = 102
$ws.Range("E43").Value = 1

# Row 44
$ws.Range("B44").Value = 95

# Row 47
$ws.Range("B47").Value = 58

# Row 56
$ws.Range("B56").Value = 14

# Row 57
$ws.Range("B57").Value = 14

# Row 58
$ws.Range("A58").Value = "Lanzarote"
$ws.Range("B58").Value = 9
$ws.Range("D58").Value = 4

# Row 59
$ws.Range("A59").Value = "Arroyo de la Luz"
$ws.Range("B59").Value = 7
$ws.Range("D59").Value = 7

# Row 60
$ws.Range("A60").Value = "Ceuta"
$ws.Range("B60").Value = 5
$ws.Range("D60").Value = 5

# Row 62
$ws.Range("B62").Value = 2
